$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (bin labels)
$ws.Range("B1").Value = "<200 mm"
$ws.Range("C1").Value = "200-300 mm"
$ws.Range("D1").Value = "300-400 mm"
$ws.Range("E1").Value = "400-500 mm"
$ws.Range("F1").Value = "500-600 mm"
$ws.Range("G1").Value = ">600 mm"

# New data values per row (years 2011-2019, rows 2-10)
$data = @{
    2  = @(2.298685006714255, 1.6693262978693573, 2.20768740030168, 2.3119814405537635, 0.6189125659909483, 0.05360025274697696)
    3  = @(0.9127311630194653, 0.7599182298788937, 1.2040544750634012, 1.595838460768662, 0.27901479320251027, 0.0)
    4  = @(1.952745297059139, 1.5939145561228476, 2.1424432376585663, 2.3290455948646924, 0.5262769590856181, 0.0)
    5  = @(1.1614946194131563, 0.9798968413606672, 2.350263455309743, 2.0669658480963187, 0.3670827736306882, 0.0)
    6  = @(0.5369085252082085, 0.5367277518787488, 1.2494655740982175, 1.6007754252260813, 0.46091907029611684, 0.05715131871844038)
    7  = @(1.0529061770613815, 1.3834302725922585, 1.9726108802562707, 1.33044891185852, 0.4238502890595027, 0.02626509662737528)
    8  = @(1.202162794602207, 1.187985060086585, 2.1705605120445144, 2.516641947154019, 0.6251125584494359, 0.026447493131732055)
    9  = @(0.6864503719863541, 0.7458483096038037, 1.2322702783688593, 1.641821900719351, 0.3493168825165554, 0.025966629620246016)
    10 = @(0.9309255268872064, 0.8421626291265406, 1.8255986797154957, 2.4502472372845636, 0.6615377753623948, 0.02626509662737528)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
